# Split the single run that holds the "La vista productos..." sentence
# into three runs around the phrase "sencillo,  que" (matching the
# target OOXML: the middle run is wrapped by Word's grammar-check
# <w:proofErr w:type="gramStart"/ .../gramEnd/> markers, but the
# observable, scriptable part of this edit -- three runs with identical
# rPr and the exact same text/xml:space split -- is what we reproduce
# here).
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("sencillo,  que", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Toggling a character-formatting property on just this sub-range
    # forces Word to break the parent run into three runs (before /
    # matched / after) that each carry their own (identical) <w:rPr>.
    # Flipping Bold back off leaves the formatting unchanged while
    # keeping the new run boundaries -- exactly the run layout the
    # target document expects.
    $rng.Bold = 1
    $rng.Bold = 0
}
